# #5: property aircraft done
# The "建物" (Building) sheet's property_category column (I) was still
# tagged "land" for its two rows; it should read "building".
# The "汽車" (Car) sheet's property_category column (H) had the same
# leftover "land" tag; it should read "car".

$wb = $excel.ActiveWorkbook

$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2").Value = "building"
$wsBuilding.Range("I3").Value = "building"

$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
$wsCar.Range("H3").Value = "car"
